$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    # Force the cell to keep a literal text value (Price column strings
    # such as "62.876.74" or "583.43" would otherwise be auto-coerced
    # into numbers by Excel). Apply a text format while writing, then
    # restore the default "Normal" style so no stray formatting/style
    # index is left behind on the cell.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# --- Update Price (D) and Volume(1h) (E) values for unchanged-identity rows ---
Set-TextValue $ws.Range("D2") "62.876.74"
$ws.Range("E2").Value = "  +2.16%  "
Set-TextValue $ws.Range("D3") "3.484.49"
$ws.Range("E3").Value = "  +2.68%  "
$ws.Range("E4").Value = "  +0.05%  "
Set-TextValue $ws.Range("D5") "583.43"
$ws.Range("E5").Value = "  +1.15%  "
Set-TextValue $ws.Range("D6") "147.44"
$ws.Range("E6").Value = "  +4.34%  "
$ws.Range("E7").Value = "  -0.11%  "
Set-TextValue $ws.Range("D8") "0.480"
$ws.Range("E8").Value = "  +1.04%  "
Set-TextValue $ws.Range("D9") "7.68"
$ws.Range("E9").Value = "  -0.45%  "
$ws.Range("E10").Value = "  +2.11%  "
Set-TextValue $ws.Range("D11") "0.400"
$ws.Range("E11").Value = "  +3.18%  "
Set-TextValue $ws.Range("D12") "4.081.54"
$ws.Range("E12").Value = "  +2.71%  "
Set-TextValue $ws.Range("D13") "30.05"
$ws.Range("E13").Value = "  +6.00%  "
Set-TextValue $ws.Range("D14") "0.126"
$ws.Range("E14").Value = "  +0.24%  "
Set-TextValue $ws.Range("D15") "3.491.11"
$ws.Range("E15").Value = "  +2.93%  "
$ws.Range("E16").Value = "  +0.66%  "
Set-TextValue $ws.Range("D17") "62.974.36"
$ws.Range("E17").Value = "  +2.29%  "
Set-TextValue $ws.Range("D18") "6.34"
$ws.Range("E18").Value = "  +3.10%  "
Set-TextValue $ws.Range("D19") "14.41"
$ws.Range("E19").Value = "  +5.17%  "
Set-TextValue $ws.Range("D20") "9.35"
$ws.Range("E20").Value = "  +4.21%  "
Set-TextValue $ws.Range("D21") "390.95"
$ws.Range("E21").Value = "  -0.04%  "
Set-TextValue $ws.Range("D22") "0.566"
$ws.Range("E22").Value = "  +2.19%  "
Set-TextValue $ws.Range("D23") "75.17"
$ws.Range("E23").Value = "  -0.47%  "
Set-TextValue $ws.Range("D24") "1.00"
$ws.Range("E24").Value = "  -0.11%  "
Set-TextValue $ws.Range("D25") "3.628.35"
Set-TextValue $ws.Range("D26") "0.0000117"
$ws.Range("E26").Value = "  +3.20%  "
Set-TextValue $ws.Range("D27") "0.180"
$ws.Range("E27").Value = "  -6.81%  "
Set-TextValue $ws.Range("D28") "7.68"
$ws.Range("E28").Value = "  +5.53%  "
$ws.Range("E29").Value = "  +0.11%  "
Set-TextValue $ws.Range("D30") "8.25"
$ws.Range("E30").Value = "  +2.60%  "
$ws.Range("E31").Value = "  -0.33%  "
$ws.Range("E32").Value = "  +3.02%  "
Set-TextValue $ws.Range("D34") "23.82"
$ws.Range("E34").Value = "  +1.81%  "
Set-TextValue $ws.Range("D39") "1.58"
$ws.Range("E39").Value = "  +7.06%  "
Set-TextValue $ws.Range("D40") "3.517.85"
$ws.Range("E40").Value = "  +2.63%  "
Set-TextValue $ws.Range("D41") "0.0771"
$ws.Range("E41").Value = "  -0.32%  "
Set-TextValue $ws.Range("D42") "0.808"
$ws.Range("E42").Value = "  +3.62%  "
Set-TextValue $ws.Range("D45") "1.72"
$ws.Range("E45").Value = "  +3.08%  "
$ws.Range("E46").Value = "  +4.17%  "
Set-TextValue $ws.Range("D47") "2.608.98"
$ws.Range("E47").Value = "  +5.69%  "
Set-TextValue $ws.Range("D48") "23.50"
$ws.Range("E48").Value = "  +1.78%  "
$ws.Range("E49").Value = "  +9.38%  "
Set-TextValue $ws.Range("D50") "6.80"
$ws.Range("E50").Value = "  +1.60%  "
$ws.Range("E51").Value = "  -0.07%  "

# --- Rows whose coin identity (and data) changed/swapped ---
$ws.Range("B35").Value = "Aptos"
$ws.Range("C35").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D35") "7.12"
$ws.Range("E35").Value = "  +2.40%  "

$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D36") "5.30"
$ws.Range("E36").Value = "  +5.03%  "

$ws.Range("B37").Value = "EnergySwap"
$ws.Range("C37").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D37") "31.64"
$ws.Range("E37").Value = "  +21.22%  "

$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D38") "171.80"
$ws.Range("E38").Value = "  +2.63%  "

$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D43") "4.47"
$ws.Range("E43").Value = "  +0.72%  "

$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D44") "42.13"
$ws.Range("E44").Value = "  -0.87%  "

